$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 317
$ws.Range("D2").Value = 305
$ws.Range("C3").Value = 53
$ws.Range("D3").Value = 48
$ws.Range("C4").Value = 1709
$ws.Range("D4").Value = 1691
$ws.Range("C5").Value = 28
$ws.Range("D5").Value = 29.5
$ws.Range("C6").Value = 237
$ws.Range("D6").Value = 230.5
$ws.Range("C7").Value = 107
$ws.Range("D7").Value = 97
$ws.Range("C8").Value = 96
$ws.Range("D8").Value = 54.5
$ws.Range("C9").Value = 54
$ws.Range("D9").Value = 49
$ws.Range("C10").Value = 242
$ws.Range("D10").Value = 231
$ws.Range("C11").Value = 98
$ws.Range("D11").Value = 56
$ws.Range("C12").Value = 49
$ws.Range("D12").Value = 41.5
$ws.Range("C13").Value = 161
$ws.Range("D13").Value = 147
$ws.Range("C14").Value = 118
$ws.Range("D14").Value = 109.5
$ws.Range("C15").Value = 37
$ws.Range("D15").Value = 34.5
$ws.Range("C16").Value = 89
$ws.Range("D16").Value = 87.5
$ws.Range("C17").Value = 33
$ws.Range("D17").Value = 33
$ws.Range("C18").Value = 140
$ws.Range("D18").Value = 116
$ws.Range("C19").Value = 125
$ws.Range("D19").Value = 120
$ws.Range("C20").Value = 146
$ws.Range("D20").Value = 136
$ws.Range("C21").Value = 46
$ws.Range("D21").Value = 27
$ws.Range("C22").Value = 41
$ws.Range("D22").Value = 25.5
$ws.Range("C23").Value = 47
$ws.Range("D23").Value = 36
$ws.Range("C24").Value = 262
$ws.Range("D24").Value = 253.5
$ws.Range("C25").Value = 57
$ws.Range("D25").Value = 51.5
$ws.Range("C26").Value = 30
$ws.Range("D26").Value = 48
$ws.Range("C27").Value = 514
$ws.Range("D27").Value = 496
$ws.Range("C28").Value = 186
